$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows describing additional NewDataSet tables (Table[2] and Table[3])
$newRows = @(
    @("/NewDataSet/Table[2]/Town", "[A-Z a-z].*"),
    @("/NewDataSet/Table[2]/County", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[2]/PostCode", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[3]/Town", "[A-Z a-z].*"),
    @("/NewDataSet/Table[3]/County", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[3]/PostCode", "[A-Z a-z 0-9].*")
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][1]
}

$ws.Range("B11").Select()
